$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The test case "Alta Deportista - Cancelar carga" (row 6) is removed from the
# test-case sheet: delete the entire row, which shifts every row below it up
# by one (rows 7/8 become 6/7) and shrinks the used range accordingly.
$ws.Rows("6").Delete()

# Conditional formats keyed off specific rows don't auto-follow a row delete
# in every case, so line them back up with the data that moved:
#  - the True/False icon rule for G6 (old row 6's own rule) goes away with it
#  - old row 7's icon rule becomes the new row 6's rule
#  - old row 8's icon rule becomes the new row 7's rule
$colG = $ws.Range("G2:G8")
$colG.FormatConditions.Item(4).Delete()
$colG.FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G7"))
$rG7 = $ws.Range("G2:G8")
$rG7.FormatConditions.Item(4).ModifyAppliesToRange($ws.Range("G6"))
$rG7.FormatConditions.Item(5).ModifyAppliesToRange($ws.Range("G7"))

$colN = $ws.Range("N2:N8")
$colN.FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("N2:N7"))

# Leave the selection where Excel puts it after deleting a row: the whole row
# that slid up into the gap.
$ws.Rows("6").Select()
